# This workbook's data rows (2-4) were cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Apply the resulting cell-level changes directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (becomes old row 4's data)
$ws.Range("A2").Value = 111697636
$ws.Range("B2").Value = 88489
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 1962
$ws.Range("F2").Value = "Vaddporing"
$ws.Range("G2").Value = "Anomoporia kamtschatica"
$ws.Range("H2").Value = "(Parmasto) Bondartseva"
$ws.Range("J2").Value = "fruktkroppar"
$ws.Range("M2").Value = $null
$ws.Range("Q2").Value = 373112.5181173298
$ws.Range("R2").Value = 6865358.590016441
$ws.Range("AC2").Value = "Växer under rötad gammal silverved"
$ws.Range("AI2").Value = "Kontinuitetsskog. Tallskog"
$ws.Range("AJ2").Value = "tall"
$ws.Range("AK2").Value = "Pinus sylvestris"
$ws.Range("AO2").Value = "Pinus sylvestris"

# Row 3 (becomes old row 2's data)
$ws.Range("A3").Value = 111697304
$ws.Range("Q3").Value = 373090.8741807578
$ws.Range("R3").Value = 6865424.499624529
$ws.Range("Z3").Value = "19:00"
$ws.Range("AB3").Value = "19:00"
$ws.Range("AI3").Value = "Luckig tallskog. K-skog"
$ws.Range("AJ3").Value = $null
$ws.Range("AK3").Value = $null
$ws.Range("AO3").Value = $null

# Row 4 (becomes old row 3's data)
$ws.Range("A4").Value = 111697236
$ws.Range("B4").Value = 8377
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 106545
$ws.Range("F4").Value = "Mindre märgborre"
$ws.Range("G4").Value = "Tomicus minor"
$ws.Range("H4").Value = "(Hartig, 1834)"
$ws.Range("J4").Value = $null
$ws.Range("M4").Value = "färska gnagspår"
$ws.Range("Q4").Value = 373121.3523494597
$ws.Range("R4").Value = 6865443.651501717
$ws.Range("Z4").Value = "00:00"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = $null
$ws.Range("AI4").Value = "Tallskog. Kontinuitetsskog"
